$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Control 47
$ws.Range("D3").Value = 0.9999999996839979
$ws.Range("E3").Value = 0.9999999996839979

# Row 4 - Control 15
$ws.Range("D4").Value = 0.9999999986968737
$ws.Range("E4").Value = 0.9999999986968737

# Row 5 - Control 7
$ws.Range("D5").Value = 0.002278140107529631
$ws.Range("E5").Value = 0.002278140107529631

# Row 6 - MDD 35
$ws.Range("D6").Value = [double]"1.82815514027753E-32"
$ws.Range("E6").Value = [double]"1.82815514027753E-32"

# Row 7 - MDD 22
$ws.Range("D7").Value = [double]"4.084263089582548E-09"
$ws.Range("E7").Value = 0.999999995915737

# Row 8 - MDD 50
$ws.Range("D8").Value = 0.9999999999999885
$ws.Range("E8").Value = [double]"1.154631945610163E-14"

# Row 9 - MDD 45
$ws.Range("D9").Value = 0.003344851357330871
$ws.Range("E9").Value = 0.9966551486426691

# Row 10 - MDD 28
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = 0.9999981597517494
$ws.Range("E10").Value = [double]"1.840248250561949E-06"

# Row 11 - Control 16
$ws.Range("D11").Value = [double]"0.0004542034276673773"
$ws.Range("E11").Value = 0.9995457965723327
$ws.Range("F11").Value = 15.93203926086426
$ws.Range("G11").Value = 0.4
